$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Text corrections (shared strings) ---
$ws.Range("F15").Value = "60000 - 80000 "
$ws.Range("F17").Value = "heutige Desktop-PCs"
$ws.Range("F16").Value = "ist normaler Wert für"

# --- New formulas replacing/augmenting static values ---
$ws.Range("D3").Formula = "=AVERAGE(A2:A6)"
$ws.Range("D4").Formula = "=VAR(A2:A6)"

$ws.Range("D10").Formula = "=AVERAGE(A9:A13)"
$ws.Range("D11").Formula = "=VAR(A9:A13)"

$ws.Range("C17").Value = "Mittelwert"
$ws.Range("D17").Formula = "=AVERAGE(A16:A20)"
$ws.Range("D18").Formula = "=VAR(A16:A20)"

# --- View / print layout changes ---
$ws.Activate()
$excel.ActiveWindow.View = "pageLayout"
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
